$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column X = 24 (day 23), Column AG = 33 (total)

# Row 2 - Bibi Cell Mundi
$ws.Cells.Item(2, 24).Value = 6645.32
$ws.Cells.Item(2, 33).Value = 198145.08

# Row 3 - Bibi Cell Vieiralves
$ws.Cells.Item(3, 24).Value = 4348
$ws.Cells.Item(3, 33).Value = 100803.01

# Row 4 - Bibi Cell Manauara
$ws.Cells.Item(4, 24).Value = 2146
$ws.Cells.Item(4, 33).Value = 68371.9

# Row 5 - Bibi Cell Ponta Negra
$ws.Cells.Item(5, 24).Value = 1171
$ws.Cells.Item(5, 33).Value = 60496.79

# Row 6 - total
$ws.Cells.Item(6, 24).Value = 14310.32
$ws.Cells.Item(6, 33).Value = 427816.78
